# Correcting Relevance Markers Walker (2018) - Wolters (2018)
# Update the metrics in row 3 (file_name = metrics_sim_with_priors.json)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.8947368421052632
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

$ws.Range("H3").Value = 0.8144859813084112
$ws.Range("I3").Value = 0.03486950184661372
$ws.Range("J3").Value = 0.7894736842105263
$ws.Range("K3").Value = 158.578947368421

$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 10
$ws.Range("S3").Value = 41
$ws.Range("T3").Value = 173
$ws.Range("U3").Value = 316
$ws.Range("V3").Value = 4258
$ws.Range("W3").Value = 4251
$ws.Range("X3").Value = 4220
$ws.Range("Y3").Value = 4088
$ws.Range("Z3").Value = 3945

$ws.Range("AF3").Value = 0.999296
$ws.Range("AG3").Value = 0.997653
$ws.Range("AH3").Value = 0.990378
$ws.Range("AI3").Value = 0.959399
$ws.Range("AJ3").Value = 0.925839
